$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "JSU(-0.8411721664047609, 1.1079255610071932, 0.8074989541855774, 2.689652300921013)"
$ws.Range("C2").Value = "JSU(-1.3025415079112728, 1.1861666495079528, 2.489566841598474, 4.753800340669389)"
$ws.Range("D2").Value = "NCT(2.6273528832259583, 1.328450087333281, -0.5675657333207895, 2.4780892838916455)"
$ws.Range("E2").Value = "NIG(1.0006319157955805, 0.8076371833616326, 4.081731228388474, 4.208562412139723)"
